# [PHOENIX-5914] CHANGES in council management
# Update the approver name for the "commissioner1" row (row 6) on the
# "approvalDetails" sheet from "Ravindra Babu/ADM_Commissioner_1" to
# "Ravindra Babu/ADM_Commissioner_2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

$ws.Range("D6").Value = "Ravindra Babu/ADM_Commissioner_2"

# Reflect the new active cell selection left behind in the sheet view.
$ws.Range("D21").Select()
